$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing D-column values that were 0 and now have small decimal values ---
# Low category (rows 69-73, years 2017-2021)
$ws.Cells.Item(69, 4).Value = 0.05161155183774097
$ws.Cells.Item(70, 4).Value = 0.08503612057176439
$ws.Cells.Item(71, 4).Value = 0.1400463483560487
$ws.Cells.Item(72, 4).Value = 0.2304786089687104
$ws.Cells.Item(73, 4).Value = 0.3788620633930405

# Medium category (rows 170-174, years 2017-2021)
$ws.Cells.Item(170, 4).Value = 0.06834730452848502
$ws.Cells.Item(171, 4).Value = 0.1126269054631187
$ws.Cells.Item(172, 4).Value = 0.1855308978543403
$ws.Cells.Item(173, 4).Value = 0.3054562194207077
$ws.Cells.Item(174, 4).Value = 0.5024408527807227

# High category (rows 271-275, years 2017-2021)
$ws.Cells.Item(271, 4).Value = 0.08952346006039959
$ws.Cells.Item(272, 4).Value = 0.1475278345532837
$ws.Cells.Item(273, 4).Value = 0.2430384466859748
$ws.Cells.Item(274, 4).Value = 0.4001768768895833
$ws.Cells.Item(275, 4).Value = 0.6583558679340189

# --- Append new "CP4All" category rows (305-405), years 1950-2050 ---

# Copy the bordered/bold style used by column A (index column) down to the new rows
$ws.Range("A2:A102").Copy($ws.Range("A305:A405"))

# Column C: category label constant across the whole new block
$ws.Range("C305:C405").Value = "CP4All"

$dValues = @(0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0.2002564859148577,0.2985841944313958,0.4450721681942594,0.6631631818984774,0.9875338019010728,1.46926280756878,2.183118297169362,3.23751779317648,4.787458843835854,7.049830103417584,10.31835166179637,14.97108536890602,21.45652596398116,30.23669067589359,41.66553897017682,55.80459339702976,72.23622117975864,90,107.7637788202414,124.1954066029702,138.3344610298232,149.7633093241064,158.5434740360188,165.028914631094,169.6816483382036,172.9501698965824,175.2125411561641,176.7624822068235,177.8168817028306,178.5307371924312,179.0124661980989,179.3368368181015,179.5549278318058,179.7014158055686)

for ($i = 0; $i -lt 101; $i++) {
  $r = 305 + $i
  $ws.Cells.Item($r, 1).Value = 303 + $i
  $ws.Cells.Item($r, 2).Value = 1950 + $i
  $ws.Cells.Item($r, 4).Value = $dValues[$i]
}
